# Apply corrected error-estimation / projected-years results to
# SoIB_summaries.xlsx, per the commit:
# "Results generated after correcting error estimation and the number
#  of projected years."

$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 5
$ws1.Range("C2").Value = 27
$ws1.Range("D2").Value = 11.6
$ws1.Range("E2").Value = 31.8

$ws1.Range("B3").Value = 7
$ws1.Range("C3").Value = 19
$ws1.Range("D3").Value = 16.3
$ws1.Range("E3").Value = 22.4

$ws1.Range("B4").Value = 10
$ws1.Range("C4").Value = 26
$ws1.Range("D4").Value = 23.3
$ws1.Range("E4").Value = 30.6

$ws1.Range("B5").Value = 5
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 11.6
$ws1.Range("E5").Value = 3.5

$ws1.Range("B6").Value = 16
$ws1.Range("C6").Value = 10
$ws1.Range("D6").Value = 37.2
$ws1.Range("E6").Value = 11.8

$ws1.Range("B7").Value = 38
$ws1.Range("C7").Value = 78

# --- Sheet: "Species qualification" -------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 43
$ws4.Range("C4").Value = 85
